$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared header string "unnamed: 1_level_1" -> "total" (cell B2)
$ws.Range("B2").Value = "total"

# Correct the off-by-one data alignment: each row label keeps its own
# data, but the two blank section-header rows (5 and 8) are removed so
# everything shifts up by the appropriate amount.
$rows = @(
    @(4, "               brasil", 1.45, 4.08, 2.08, 2.85, 2.93, 2.05, 3.3, 3.37),
    @(5, "urbana", 1.53, 4.46, 2.18, 2.96, 3.04, 2.09, 3.33, 3.4),
    @(6, "rural", 6.03, 9.11, 6.39, 13.04, 15.67, 11.77, 27.64, 22.69),
    @(7, "norte", 4.97, 12.89, 6.57, 7.21, 6.6, 5.83, 8.44, 8.4),
    @(8, "rondônia", 9.23, 19.25, 10.94, 16.23, 18.75, 14.58, 18.56, 24.88),
    @(9, "acre", 16.39, 31.76, 20.73, 27.55, 20.87, 20.51, 34.54, 23.58),
    @(10, "amazonas", 11.67, 26.02, 15.53, 23.65, 24.33, 14.66, 23.01, 32.37),
    @(11, "roraima", 11.46, 39.87, 19.58, 22.42, 27.7, 11.26, 37.31, 31.37),
    @(12, "pará", 7.69, 21.87, 9.75, 9.98, 8.27, 8.84, 13.64, 12.95),
    @(13, "amapá", 16.65, 89.69, 20.44, 24.55, 24.13, 19.71, 23.84, 34.71),
    @(14, "tocantins", 13.44, 30.83, 15.07, 31.69, 30.1, 17.93, 27.37, 18.87),
    @(15, "nordeste", 2.72, 6.36, 3.82, 5.5, 5.46, 3.72, 6.39, 6.25),
    @(16, "maranhão", 15.26, 25.54, 19.39, 27.46, 28.16, 16.37, 34.27, 23.76),
    @(17, "piauí", 14.58, 16.4, 18.86, 24.29, 30.07, 26.81, 31.83, 30.96),
    @(18, "ceará", 4.62, 15.94, 6.9, 6.95, 9.27, 6.52, 12.58, 11.77),
    @(19, "rio grande do norte", 7.33, 21.63, 13.04, 19.12, 17.04, 11.0, 20.8, 31.39),
    @(20, "paraíba", 9.75, 22.21, 14.28, 22.53, 18.26, 18.1, 35.37, 30.99),
    @(21, "pernambuco", 5.39, 14.42, 9.42, 11.26, 12.46, 6.75, 13.31, 12.23),
    @(22, "alagoas", 13.17, 26.38, 18.23, 23.75, 31.26, 21.92, 37.12, 30.04),
    @(23, "sergipe", 7.52, 26.87, 12.89, 22.48, 17.11, 17.18, 18.04, 24.01),
    @(24, "bahia", 5.27, 12.38, 6.64, 10.09, 9.69, 7.59, 11.93, 11.15),
    @(25, "sudeste", 2.53, 8.01, 3.79, 5.0, 5.69, 3.58, 5.56, 5.29),
    @(26, "minas gerais", 4.66, 13.73, 6.94, 9.35, 9.8, 7.16, 11.32, 10.63),
    @(27, "espírito santo", 8.94, 24.23, 14.32, 16.57, 20.65, 14.72, 24.66, 25.0),
    @(28, "rio de janeiro", 4.48, 17.41, 7.4, 9.13, 10.93, 6.57, 10.57, 8.8),
    @(29, "são paulo", 3.77, 11.64, 5.7, 7.52, 8.49, 5.17, 8.03, 7.38),
    @(30, "sul", 3.32, 11.73, 5.13, 6.59, 7.19, 4.99, 7.26, 7.48),
    @(31, "paraná", 4.5, 17.69, 7.46, 9.93, 11.51, 7.1, 12.5, 11.15),
    @(32, "santa catarina", 11.47, 19.79, 14.97, 20.57, 20.65, 17.84, 24.94, 19.93),
    @(33, "rio grande do sul", 4.74, 23.37, 7.34, 9.42, 9.47, 6.49, 8.92, 10.06),
    @(34, "centro-oeste", 3.39, 10.01, 5.01, 7.34, 7.0, 5.16, 9.13, 7.65),
    @(35, "mato grosso do sul", 15.62, 32.93, 20.41, 23.03, 24.82, 15.8, 27.31, 27.8),
    @(36, "mato grosso", 6.64, 29.52, 10.0, 21.02, 17.06, 11.01, 19.27, 17.57),
    @(37, "goiás", 5.08, 11.94, 7.17, 10.26, 10.47, 8.85, 15.16, 12.13),
    @(38, "distrito federal", 5.54, 28.74, 8.72, 13.4, 11.42, 7.94, 14.3, 12.82)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}

# Rows 39 and 40 are now duplicates of the shifted-up data; remove them
# so the sheet dimension shrinks from A1:I40 to A1:I38.
$ws.Rows("39:40").Delete()
